# Add updateDB method to TC005 Completed TC008, RC009
# Populate the "Sponsor" and "Client" sheets with their DB id lists and
# restore the view state (active sheet / selected cell) to what it was
# after the update: Sponsor becomes the active tab with D14 selected,
# while Client keeps its J14 selection but is no longer the active tab.

$wb = $excel.ActiveWorkbook

$wsSponsor = $wb.Worksheets.Item("Sponsor")
$wsClient  = $wb.Worksheets.Item("Client")

# Sponsor (sheet1) data
$wsSponsor.Range("A1").Value = 2948492
$wsSponsor.Range("A2").Value = 2948532
$wsSponsor.Range("A3").Value = 2948534

# Client (sheet2) data
$wsClient.Range("A1").Value = 2948512
$wsClient.Range("A2").Value = 2948533
$wsClient.Range("A3").Value = 2948535

# Restore selections on each sheet. Select Client first so that the final
# active/selected sheet ends up being Sponsor, matching the target state.
$wsClient.Range("J14").Select()
$wsSponsor.Range("D14").Select()
